$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing hyperlinks so we can rebuild them cleanly at the new row positions
$ws.Cells.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2025-09-25 18:27:05'
$ws.Range("B2").Value = '【SES案件多数】バックエンドエンジニア募集(Java/PHP/Python/Node.js)'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5399874'
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5399874')
$ws.Range("G2").Value = 320
$ws.Range("H2").Value = '🔥Python ★Java ◆Node.js ○PHP'

# Row 3
$ws.Range("A3").Value = '2025-09-25 18:27:05'
$ws.Range("B3").Value = '【低予算希望】LINE公式アカウント+社食注文システム開発依頼(社内利用のみ)'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5400375'
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5400375')
$ws.Range("G3").Value = 118
$ws.Range("H3").Value = '◆開発,システム開発'

# Row 4
$ws.Range("A4").Value = '2025-09-25 18:27:05'
$ws.Range("B4").Value = '【急募】住宅展示場マッチング診断サービスのMVP開発依頼'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5399759'
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5399759')
$ws.Range("G4").Value = 75
$ws.Range("H4").Value = '◆開発'

# Row 5
$ws.Range("A5").Value = '2025-09-25 18:27:05'
$ws.Range("B5").Value = '【フリーランス募集】CTビューアーソフト気道抽出機能開発'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5400101'
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5400101')
$ws.Range("G5").Value = 68
$ws.Range("H5").Value = '◆開発'

# Row 6
$ws.Range("A6").Value = '2025-09-25 18:27:05'
$ws.Range("B6").Value = 'MYSQLからGoogleスプレッドシートへデータ取り込み及びスプレッドシート改修'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5400606'
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5400606')
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = '◇MySQL'

# Row 7
$ws.Range("A7").Value = '2025-09-25 18:27:05'
$ws.Range("B7").Value = '【急募】SOLIDWORKS2024での機械設計と製図依頼'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5400338'
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5400338')
$ws.Range("G7").Value = 25

# Row 8
$ws.Range("A8").Value = '2025-09-25 18:27:05'
$ws.Range("B8").Value = '当社CTソフトへの機能追加:気道抽出'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '3,000,000 円 ~ 5,000,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5400094'
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5400094')
$ws.Range("G8").Value = 25

# Row 9
$ws.Range("A9").Value = '2025-09-25 18:27:05'
$ws.Range("B9").Value = '【SES案件多数/リモート可】インフラエンジニア募集(AWS/Linux/NW設計・構築 等歓迎)'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5399876'
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5399876')
$ws.Range("G9").Value = 25

# Row 10
$ws.Range("A10").Value = '2025-09-25 18:27:05'
$ws.Range("B10").Value = '【SES案件多数/リモート可】フルスタックエンジニア募集(フロント〜バック〜クラウドまで歓迎)'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5399877'
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5399877')
$ws.Range("G10").Value = 25

# Row 11
$ws.Range("A11").Value = '2025-09-25 18:27:05'
$ws.Range("B11").Value = '【急募】Nuxt3でのWEBページ表示速度改善依頼'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5400231'
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5400231')
$ws.Range("G11").Value = 18

# Row 12
$ws.Range("A12").Value = '2025-09-25 18:27:05'
$ws.Range("B12").Value = '限定公開 PR 限定公開の仕事'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5399347'
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5399347')
$ws.Range("G12").Value = 13

# Row 13
$ws.Range("A13").Value = '2025-09-25 18:27:05'
$ws.Range("B13").Value = '【急募】iPhoneのwifiMACアドレスを偽装する方法を教えてください'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5400676'
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5400676')
$ws.Range("G13").Value = 10

# Row 14
$ws.Range("A14").Value = '2025-09-25 18:27:05'
$ws.Range("B14").Value = '【急募】ライフプランシミュレーターのバグ確認と使用感調査'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '~ 5,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5400626'
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5400626')
$ws.Range("G14").Value = 10

# Row 15
$ws.Range("A15").Value = '2025-09-25 18:27:05'
$ws.Range("B15").Value = '【SalesIQ活用】CRMと連携したリード獲得方法を教えてください'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '~ 5,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5400402'
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5400402')
$ws.Range("G15").Value = 10
